$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.011.48"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.900.94"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7412"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.65"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3065"
$ws.Range("E8").Value = "  -3.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.87"
$ws.Range("E9").Value = "  -6.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06905"
$ws.Range("E10").Value = "  -2.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08003"
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7602"
$ws.Range("E12").Value = "  -2.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.897.88"
$ws.Range("E13").Value = "  -1.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.235"
$ws.Range("E14").Value = "  -2.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.31"
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.213"
$ws.Range("E16").Value = "  +2.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.031.59"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.03"
$ws.Range("E18").Value = "  -3.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007758"
$ws.Range("E19").Value = "  -2.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.33"
$ws.Range("E20").Value = "  -5.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.152.07"
$ws.Range("E22").Value = "  -1.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9999"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.075"
$ws.Range("E24").Value = "  +5.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.305"
$ws.Range("E25").Value = "  -2.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.38"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.81"
$ws.Range("E27").Value = "  -1.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1255"
$ws.Range("E28").Value = "  -3.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.038"
$ws.Range("E29").Value = "  -6.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.351"
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.535"
$ws.Range("E31").Value = "  -1.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.299"
$ws.Range("E32").Value = "  -2.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.045"
$ws.Range("E33").Value = "  -2.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05279"
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.292"
$ws.Range("E35").Value = "  -1.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7399"
$ws.Range("E36").Value = "  -2.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.726"
$ws.Range("E37").Value = "  -1.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01939"
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.778"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.272"
$ws.Range("E40").Value = "  -3.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4445"
$ws.Range("E41").Value = "  -1.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.09"
$ws.Range("E42").Value = "  -6.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.965"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.0000"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8370"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.631"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.27"
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.816"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.051.07"
$ws.Range("E49").Value = "  -1.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.73"
$ws.Range("E50").Value = "  -3.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1169"
$ws.Range("E51").Value = "  -5.29%  "
